$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1043.3846
$ws.Range("I28").Value = 575.05554
$ws.Range("J28").Value = 2097.125
$ws.Range("K28").Value = 575.05554
$ws.Range("L28").Value = 2097.125
$ws.Range("M28").Value = -90.05553999999995
$ws.Range("N28").Value = -3067.125
$ws.Range("H40").Value = 3332
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 1619.75
$ws.Range("I41").Value = 2245.5
$ws.Range("J41").Value = 994
$ws.Range("K41").Value = 2245.5
$ws.Range("L41").Value = 994
$ws.Range("M41").Value = -1805.5
$ws.Range("N41").Value = -1874
$ws.Range("H106").Value = 4185.3076
$ws.Range("I106").Value = 4426.5
$ws.Range("K106").Value = 4426.5
$ws.Range("M106").Value = -3795.5
$ws.Range("H112").Value = 2610.7778
$ws.Range("J112").Value = 2828.2856
$ws.Range("L112").Value = 8484.856800000001
$ws.Range("N112").Value = -10700.8568
$ws.Range("H129").Value = 3271.3076
$ws.Range("I129").Value = 4838.4
$ws.Range("K129").Value = 14515.2
$ws.Range("M129").Value = -9515.199999999999
$ws.Range("H137").Value = 1669.8422
$ws.Range("I137").Value = 1557
$ws.Range("K137").Value = 4671
$ws.Range("M137").Value = -2121
$ws.Range("H138").Value = 3521.625
$ws.Range("I138").Value = 2672.4
$ws.Range("J138").Value = 4270.9414
$ws.Range("K138").Value = 8017.200000000001
$ws.Range("L138").Value = 12812.8242
$ws.Range("M138").Value = -2877.200000000001
$ws.Range("N138").Value = -23092.8242

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5165.316
$ws.Range("I32").Value = 5243.647
$ws.Range("K32").Value = 5243.647
$ws.Range("M32").Value = -4956.647
$ws.Range("H61").Value = 3494.675
$ws.Range("I61").Value = 3620.4482
$ws.Range("K61").Value = 3620.4482
$ws.Range("M61").Value = -3408.4482
$ws.Range("H122").Value = 4434
$ws.Range("I122").Value = 4220.1665
$ws.Range("K122").Value = 12660.4995
$ws.Range("M122").Value = -10210.4995
$ws.Range("H132").Value = 29418508
$ws.Range("I132").Value = 7241.154
$ws.Range("J132").Value = 125005130
$ws.Range("K132").Value = 21723.462
$ws.Range("L132").Value = 375015390
$ws.Range("M132").Value = -19193.462
$ws.Range("N132").Value = -375020450
$ws.Range("H136").Value = 3494.675
$ws.Range("I136").Value = 3620.4482
$ws.Range("K136").Value = 10861.3446
$ws.Range("M136").Value = -8311.3446

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9999
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 5739.4
$ws.Range("J8").Value = 5739.4
$ws.Range("L8").Value = 5739.4
$ws.Range("N8").Value = -6019.4
$ws.Range("H25").Value = 7374.75
$ws.Range("J25").Value = 8166.6665
$ws.Range("L25").Value = 8166.6665
$ws.Range("N25").Value = -8514.666499999999
$ws.Range("H31").Value = 1540.1212
$ws.Range("I31").Value = 1893.25
$ws.Range("K31").Value = 1893.25
$ws.Range("M31").Value = -1598.25
$ws.Range("H34").Value = 1540.1212
$ws.Range("I34").Value = 1893.25
$ws.Range("K34").Value = 1893.25
$ws.Range("M34").Value = -1691.25
$ws.Range("H105").Value = 1334
$ws.Range("I105").Value = 1410.125
$ws.Range("K105").Value = 1410.125
$ws.Range("M105").Value = 336.875

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 110.47059
$ws.Range("I40").Value = 124.96154
$ws.Range("K40").Value = 499.84616
$ws.Range("M40").Value = -430.84616
$ws.Range("H68").Value = 1675.2593
$ws.Range("J68").Value = 1772.2084
$ws.Range("L68").Value = 5316.6252
$ws.Range("N68").Value = -6938.6252
$ws.Range("H71").Value = 1675.2593
$ws.Range("J71").Value = 1772.2084
$ws.Range("L71").Value = 15949.8756
$ws.Range("N71").Value = -24061.8756
$ws.Range("H141").Value = 8516
$ws.Range("I141").Value = 8999
$ws.Range("K141").Value = 26997
$ws.Range("M141").Value = -21817

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 17516.5
$ws.Range("J52").Value = 16033
$ws.Range("L52").Value = 16033
$ws.Range("N52").Value = -16551
$ws.Range("H70").Value = 9561
$ws.Range("I70").Value = 8748.5
$ws.Range("J70").Value = 11998.5
$ws.Range("K70").Value = 8748.5
$ws.Range("L70").Value = 11998.5
$ws.Range("M70").Value = -8478.5
$ws.Range("N70").Value = -12538.5
$ws.Range("H73").Value = 9561
$ws.Range("I73").Value = 8748.5
$ws.Range("J73").Value = 11998.5
$ws.Range("K73").Value = 8748.5
$ws.Range("L73").Value = 11998.5
$ws.Range("M73").Value = -7812.5
$ws.Range("N73").Value = -13870.5
$ws.Range("H113").Value = 2988.5386
$ws.Range("I113").Value = 2988.5386
$ws.Range("K113").Value = 2988.5386
$ws.Range("M113").Value = -818.5385999999999
$ws.Range("H132").Value = 3895.8262
$ws.Range("I132").Value = 3589.5
$ws.Range("K132").Value = 10768.5
$ws.Range("M132").Value = -8238.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3182.4375
$ws.Range("I22").Value = 1360.5
$ws.Range("J22").Value = 4275.6
$ws.Range("K22").Value = 1360.5
$ws.Range("L22").Value = 4275.6
$ws.Range("M22").Value = -1065.5
$ws.Range("N22").Value = -4865.6
$ws.Range("H27").Value = 3182.4375
$ws.Range("I27").Value = 1360.5
$ws.Range("J27").Value = 4275.6
$ws.Range("K27").Value = 1360.5
$ws.Range("L27").Value = 4275.6
$ws.Range("M27").Value = -1253.5
$ws.Range("N27").Value = -4489.6
$ws.Range("H46").Value = 2522.923
$ws.Range("I46").Value = 2532.6667
$ws.Range("J46").Value = 2514.5715
$ws.Range("K46").Value = 2532.6667
$ws.Range("L46").Value = 2514.5715
$ws.Range("M46").Value = -2344.6667
$ws.Range("N46").Value = -2890.5715
$ws.Range("H136").Value = 142861630
$ws.Range("I136").Value = 5199.3335
$ws.Range("K136").Value = 15598.0005
$ws.Range("M136").Value = -13048.0005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 96999.5
$ws.Range("J119").Value = 96999.5
$ws.Range("L119").Value = 96999.5
$ws.Range("N119").Value = -106675.5
$ws.Range("H132").Value = 3083.75
$ws.Range("I132").Value = 2980.3845
$ws.Range("J132").Value = 3531.6667
$ws.Range("K132").Value = 8941.1535
$ws.Range("L132").Value = 10595.0001
$ws.Range("M132").Value = -6411.1535
$ws.Range("N132").Value = -15655.0001
